# Update "Horarios actualizados Línea 141 - 555"
# Refresh the scrape timestamp + arrival data across the three sheets of the
# workbook (LP1912, LP1912-215, 6203-6173).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

$newScrapTime = "04:29:16"

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1.Range("A2").Value = "Última actualización: " + $newScrapTime
$ws1.Range("A3").Value = "Total filas: 9"

$ws1.Range("A6").Value = $newScrapTime
$ws1.Range("B6").Value = "04:48"
$ws1.Range("C6").Value = "81_EL PELIGRO"
$ws1.Range("D6").Value = 19
$ws1.Range("E6").Value = "LP1912"

$ws1.Range("A7").Value = $newScrapTime
$ws1.Range("B7").Value = "04:53"
$ws1.Range("C7").Value = "11_ETCHEVERRY"
$ws1.Range("D7").Value = 24
$ws1.Range("E7").Value = "LP1912"

$ws1.Range("A8").Value = $newScrapTime
$ws1.Range("B8").Value = "05:17"
$ws1.Range("C8").Value = "17_ROMERO"
$ws1.Range("D8").Value = 48
$ws1.Range("E8").Value = "LP1912"

$ws1.Range("A9").Value = $newScrapTime
$ws1.Range("B9").Value = "05:22"
$ws1.Range("C9").Value = "23_HERNANDEZ"
$ws1.Range("D9").Value = 53
$ws1.Range("E9").Value = "LP1912"

$ws1.Range("A10").Value = $newScrapTime
$ws1.Range("B10").Value = "05:42"
$ws1.Range("C10").Value = "14_ABASTO"
$ws1.Range("D10").Value = 73
$ws1.Range("E10").Value = "LP1912"

$ws1.Range("A11").Value = $newScrapTime
$ws1.Range("B11").Value = "05:47"
$ws1.Range("C11").Value = "17_ROMERO"
$ws1.Range("D11").Value = 78
$ws1.Range("E11").Value = "LP1912"

$ws1.Range("A12").Value = $newScrapTime
$ws1.Range("B12").Value = "06:01"
$ws1.Range("C12").Value = "16_SANTA ANA"
$ws1.Range("D12").Value = 92
$ws1.Range("E12").Value = "LP1912"

$ws1.Range("A13").Value = $newScrapTime
$ws1.Range("B13").Value = "06:09"
$ws1.Range("C13").Value = "10_OLMOS"
$ws1.Range("D13").Value = 100
$ws1.Range("E13").Value = "LP1912"

$ws1.Range("A14").Value = $newScrapTime
$ws1.Range("B14").Value = "06:16"
$ws1.Range("C14").Value = "215A_EL PATO"
$ws1.Range("D14").Value = 107
$ws1.Range("E14").Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2.Range("A2").Value = "Última actualización: " + $newScrapTime
$ws2.Range("A3").Value = "Total filas: 1"

$ws2.Range("A5").Value = "Hora_Scrap"
$ws2.Range("B5").Value = "Hora_Llegada"
$ws2.Range("C5").Value = "Linea"
$ws2.Range("D5").Value = "Minutos"
$ws2.Range("E5").Value = "Parada"
$ws1.Range("A5:E5").Copy()
$ws2.Range("A5:E5").PasteSpecial(-4122)

$ws2.Range("A6").Value = $newScrapTime
$ws2.Range("B6").Value = "06:16"
$ws2.Range("C6").Value = "215A_EL PATO"
$ws2.Range("D6").Value = 107
$ws2.Range("E6").Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3.Range("A2").Value = "Última actualización: " + $newScrapTime
